$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1986531986531987
$ws.Range("C2").Value = 0.5387205387205387
$ws.Range("J2").Value = 0.01346801346801347
$ws.Range("P2").Value = 0.1616161616161616
$ws.Range("S2").Value = 0.08754208754208755
$ws.Range("B3").Value = 0.0119047619047619
$ws.Range("C3").Value = 0.04166666666666666
$ws.Range("J3").Value = 0.04761904761904762
$ws.Range("P3").Value = 0.6845238095238095
$ws.Range("S3").Value = 0.2142857142857143
$ws.Range("P4").Value = 0.75
$ws.Range("S4").Value = 0.25
$ws.Range("B6").Value = 0.05957446808510639
$ws.Range("D6").Value = 0.00425531914893617
$ws.Range("F6").Value = 0.07234042553191489
$ws.Range("J6").Value = 0.2595744680851064
$ws.Range("O6").Value = 0.008510638297872341
$ws.Range("Q6").Value = 0.1148936170212766
$ws.Range("R6").Value = 0.06808510638297872
$ws.Range("S6").Value = 0.4127659574468085
$ws.Range("B7").Value = 0.08056872037914692
$ws.Range("D7").Value = 0.009478672985781991
$ws.Range("E7").Value = 0.009478672985781991
$ws.Range("F7").Value = 0.09004739336492891
$ws.Range("J7").Value = 0.1137440758293839
$ws.Range("O7").Value = 0.04265402843601896
$ws.Range("Q7").Value = 0.1279620853080569
$ws.Range("R7").Value = 0.08056872037914692
$ws.Range("S7").Value = 0.4454976303317535
$ws.Range("B8").Value = 0.1106870229007634
$ws.Range("D8").Value = 0.01717557251908397
$ws.Range("F8").Value = 0.06297709923664122
$ws.Range("J8").Value = 0.1049618320610687
$ws.Range("O8").Value = 0.04389312977099236
$ws.Range("Q8").Value = 0.1106870229007634
$ws.Range("R8").Value = 0.09923664122137404
$ws.Range("S8").Value = 0.4503816793893129
$ws.Range("B9").Value = 0.1142857142857143
$ws.Range("F9").Value = 0.05142857142857143
$ws.Range("J9").Value = 0.08
$ws.Range("O9").Value = 0.02285714285714286
$ws.Range("Q9").Value = 0.1657142857142857
$ws.Range("R9").Value = 0.1028571428571429
$ws.Range("S9").Value = 0.4628571428571429
$ws.Range("B10").Value = 0.1045150501672241
$ws.Range("D10").Value = 0.0117056856187291
$ws.Range("E10").Value = 0.0008361204013377926
$ws.Range("F10").Value = 0.07692307692307693
$ws.Range("J10").Value = 0.137123745819398
$ws.Range("O10").Value = 0.02424749163879599
$ws.Range("Q10").Value = 0.1622073578595318
$ws.Range("R10").Value = 0.08361204013377926
$ws.Range("S10").Value = 0.3988294314381271
$ws.Range("G11").Value = 0.1467065868263473
$ws.Range("J11").Value = 0.09580838323353294
$ws.Range("K11").Value = 0.218562874251497
$ws.Range("L11").Value = 0.5179640718562875
$ws.Range("S11").Value = 0.02095808383233533
$ws.Range("G12").Value = 0.6885245901639344
$ws.Range("J12").Value = 0.1967213114754098
$ws.Range("K12").Value = 0.01639344262295082
$ws.Range("L12").Value = 0.04918032786885246
$ws.Range("S12").Value = 0.04918032786885246
$ws.Range("G13").Value = 0.75
$ws.Range("J13").Value = 0.2
$ws.Range("S13").Value = 0.05
$ws.Range("F15").Value = 0.004149377593360996
$ws.Range("H15").Value = 0.1452282157676349
$ws.Range("I15").Value = 0.05809128630705394
$ws.Range("J15").Value = 0.2655601659751037
$ws.Range("K15").Value = 0.05394190871369295
$ws.Range("M15").Value = 0.004149377593360996
$ws.Range("O15").Value = 0.0912863070539419
$ws.Range("S15").Value = 0.3775933609958506
$ws.Range("F16").Value = 0.01104972375690608
$ws.Range("H16").Value = 0.1602209944751381
$ws.Range("I16").Value = 0.0718232044198895
$ws.Range("J16").Value = 0.3646408839779006
$ws.Range("K16").Value = 0.1270718232044199
$ws.Range("M16").Value = 0.04419889502762431
$ws.Range("O16").Value = 0.06077348066298342
$ws.Range("S16").Value = 0.1602209944751381
$ws.Range("F17").Value = 0.01208459214501511
$ws.Range("H17").Value = 0.1933534743202417
$ws.Range("I17").Value = 0.0634441087613293
$ws.Range("J17").Value = 0.4138972809667674
$ws.Range("K17").Value = 0.1178247734138973
$ws.Range("M17").Value = 0.01812688821752266
$ws.Range("O17").Value = 0.08157099697885196
$ws.Range("S17").Value = 0.09969788519637462
$ws.Range("F18").Value = 0.02487562189054726
$ws.Range("H18").Value = 0.2039800995024875
$ws.Range("I18").Value = 0.06965174129353234
$ws.Range("J18").Value = 0.3930348258706468
$ws.Range("K18").Value = 0.1144278606965174
$ws.Range("M18").Value = 0.01492537313432836
$ws.Range("O18").Value = 0.05472636815920398
$ws.Range("S18").Value = 0.1243781094527363
$ws.Range("F19").Value = 0.01990049751243781
$ws.Range("H19").Value = 0.2487562189054726
$ws.Range("I19").Value = 0.07960199004975124
$ws.Range("J19").Value = 0.3326226012793177
$ws.Range("K19").Value = 0.1115849324804549
$ws.Range("M19").Value = 0.03198294243070363
$ws.Range("N19").Value = 0.001421464108031272
$ws.Range("O19").Value = 0.05330490405117271
$ws.Range("S19").Value = 0.1208244491826581
